$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that held the "Docentes responsáveis:" value (old row 13, containing
# only "5840535 - Messias Borges Silva" in B/C) is removed entirely; everything
# below shifts up by one row.
$ws.Rows.Item(13).Delete()

# After the shift, several content cells (B/C columns) end up holding
# different text than a plain shift would produce, so set them explicitly.
$ws.Range("B10:C10").Value = "5840535 - Messias Borges Silva"
$ws.Range("B13:C13").Value = "Semestral"

# Row 15 needs the literal text "01/01/2018" (not an Excel date serial).
# Assigning via .Value would auto-convert it to a date, so instead copy the
# already-text cell B8 (which holds the same literal string) and paste its
# value into B15/C15 — this keeps it as text and preserves the row's existing
# number format / style.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B18:C18").Value = "5840535 - Messias Borges Silva"
$ws.Range("B19:C19").Value = "2 provas escritas"
$ws.Range("B20:C20").Value = "Serão avaliados os conteúdos discutidos em sala e constantes da ementa do curso. MF = (0,40*P1 + 0,40*P2 + 0,20*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários."
$ws.Range("B21:C21").Value = "Uma provas escrita com conteúdo de todo o semestre. NF = (MF + PR)/2, onde PR é uma prova de recuperação"
